$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 'The applicant has strong skills in Python, Django, ReactJS, and JavaScript, showcasing proficiency in web development and backend technologies. However, there is a lack of experience in specific technologies such as Flutter, Dart, Firestore, Firebase Authentication, and mobile app architecture/design, which are crucial for the job role. Additionally, the applicant''s experience in mobile app development is limited, and there is no mention of working under a Scrum Master, which is a requirement specified in the job description.'

$ws.Range("E2").Value = 'The applicant has demonstrated excellent communication skills, adaptability, and a strong interest in learning from Japanese work culture. The responses reflect a positive attitude, a willingness to learn and improve, and an ability to overcome challenges. The applicant''s keen interest in Japanese values and cultural alignment indicates a high level of adaptability, making them suitable for the company''s requirements. Overall, the applicant''s personality aligns with the company''s soft skill requirements, warranting a score of 4.'

$ws.Range("C3").Value = 'The applicant has a strong background and experience in relevant technologies such as Flutter, Dart, Firebase, WebRTC, NodeJS, ExpressJS, and HTML/CSS/JS. The projects indicate proficiency in developing mobile and web applications, utilizing real-time communication, and smart home device integration. However, the applicant lacks experience with Android/iOS, Firestore, Firebase Authentication, and Cloud Storage/Messaging, which are essential for the specified role in application development. Despite this, the applicant''s demonstrated knowledge and experience make them a strong candidate for the position.'

$ws.Range("E3").Value = 'The applicant exhibits a high level of adaptability, problem-solving skills, and a strong willingness to work in Japan. Their expressed interest in learning a new language and embracing new cultures, combined with their collaborative work style, indicates a positive fit for the company''s soft skill requirements.'

$ws.Range("C4").Value = 'The applicant has a strong background in computer vision, image processing, backend development, and web development. Although the applicant''s projects do not directly align with the specific technologies required by the company, the candidate possesses valuable skills such as Python, PyTorch, Tensorflow, Django Rest Framework, and ReactJS, which demonstrate a proficiency in software development. However, the lack of experience in Flutter, Dart, mobile app architecture, design patterns, Firestore, Firebase Authentication, Cloud Storage, and Cloud Messaging may be a gap in the candidate''s skill set for the specific role. The applicant''s experience in developing a library management system and backend integration with machine learning models showcases relevant experience in software development, which could be leveraged in the application development role, albeit with some additional training in mobile app development technologies.'

$ws.Range("E4").Value = 'Based on the provided answers, the applicant demonstrates strong adaptability, communication, and teamwork skills, which are essential for integrating into a new work culture in Japan. The candidate expresses a deep interest in the work culture and technological advancements in Japan, showcasing a willingness to embrace new experiences. Additionally, the applicant''s emphasis on working well in a team and readiness to learn new skills align with the company''s soft skill requirements, indicating a strong potential for success in a collaborative work environment. The candidate''s willingness to learn a new language and adapt to a new culture further solidifies their openness to the opportunity of working in Japan, thus warranting a score of 4.'

$ws.Range("C5").Value = 'The applicant''s projects demonstrate strong proficiency in front-end development using JavaScript, HTML, and CSS, and the ability to work on group projects. However, the lack of experience with Flutter, Dart, Android/iOS, Firestore, Firebase Authentication, Cloud Storage/Messaging, and mobile app architecture and design patterns, which are crucial for the job, makes the applicant less suited for the position. The candidate''s project experience does not align closely with the specific technologies and skills required by the company, resulting in a moderate score.'

$ws.Range("E5").Value = 'The applicant has displayed good communication skills, teamwork, problem-solving ability, and adaptability to new cultures, as evidenced by the responses provided. Their openness to working in Japan, the admiration for Japan''s work culture, and the interest in learning a new language demonstrate a strong willingness to work in a different environment, aligning well with the company''s soft skill requirements. As such, the applicant''s strong interpersonal skills and adaptability warrant a score of 4.'
